$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets(1))
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 1061
$ws1.Range("G3").Value = 68
$ws1.Range("F4").Value = 531
$ws1.Range("F5").Value = 13961
$ws1.Range("F6").Value = 47
$ws1.Range("F7").Value = 567
$ws1.Range("F8").Value = 222
$ws1.Range("F9").Value = 1798
$ws1.Range("F10").Value = 173
$ws1.Range("F11").Value = 145
$ws1.Range("F14").Value = 544
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 6
$ws1.Range("F18").Value = 14033
$ws1.Range("F19").Value = 372
$ws1.Range("F20").Value = 634
$ws1.Range("F21").Value = 15001
$ws1.Range("F23").Value = 8314
$ws1.Range("F26").Value = 31
$ws1.Range("F33").Value = 31
$ws1.Range("F42").Value = 222
$ws1.Range("F45").Value = 5117

# Sheet "全部类型" (Worksheets(4))
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1061
$ws4.Range("G3").Value = 68
$ws4.Range("F4").Value = 531
$ws4.Range("F5").Value = 13961
$ws4.Range("F6").Value = 47
$ws4.Range("F7").Value = 567
$ws4.Range("F8").Value = 222
$ws4.Range("F9").Value = 1798
$ws4.Range("F10").Value = 173
$ws4.Range("F11").Value = 145
$ws4.Range("F14").Value = 544
$ws4.Range("F15").Value = 35
$ws4.Range("F16").Value = 6
$ws4.Range("F18").Value = 14033
$ws4.Range("F19").Value = 372
$ws4.Range("F20").Value = 634
$ws4.Range("F21").Value = 15001
$ws4.Range("F23").Value = 8314
$ws4.Range("F26").Value = 31
$ws4.Range("F33").Value = 31
$ws4.Range("F44").Value = 222
$ws4.Range("F47").Value = 5117
